# Actors.xlsx — "Planning for another big cleanup/restructuring... Q.Q"
#
# The roster's last row (row 10 of the XML-mapped Table2: Sophie Q. / Voice)
# is being swapped out for a new actor: Huilian Q., wielding a Microphone.
# Only the Name and Weapon columns change for that row; Health/Attack/
# Defense/Speed stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Huilian Q."
$ws.Range("F10").Value = "Microphone"

# Leave the cursor where the author left it when saving.
$ws.Range("F11").Select()
